$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix latitude (G) and longitude (H) for rows 3 and 4 so they match row 2
# (the three rows are the same registration center in different languages).
$ws.Range("G3").Value = 34.52117
$ws.Range("G4").Value = 34.52117

# Copy H2's text value (with its leading non-breaking space) into H3/H4 so the
# cells stay text instead of being re-parsed as numbers.
$ws.Range("H2").Copy($ws.Range("H3"))
$ws.Range("H2").Copy($ws.Range("H4"))

# Bump number_of_kiosks (column L): the first center (rows 2-4) now has 3
# kiosks, every other center (rows 5-46) now has 2 kiosks.
for ($r = 2; $r -le 4; $r++) {
    $ws.Cells.Item($r, 12).Value = 3
}
for ($r = 5; $r -le 46; $r++) {
    $ws.Cells.Item($r, 12).Value = 2
}

# Update the active selection left by the author.
$ws.Range("L28").Select() | Out-Null
